$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename id_giria -> id_phrasal (column A1)
$ws.Range("A1").Value = "id_phrasal"

# Data row 2: replace the "giria" example entry with the new "phrasal verb" entry
$ws.Range("B2").Value = "Give up"
$ws.Range("C2").Value = "Desistir."
$ws.Range("D2").Value = "Desistir."
$ws.Range("E2").Value = "I want to Give Up."
$ws.Range("F2").Value = "Eu quero desistir."
# A2 (id value "1") and G2 (dificuldade "2") remain unchanged

# Update the view: scroll/select so column E is left-most visible and F2 is the active cell
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
[void]$ws.Range("F2").Select()
